$wb = $excel.ActiveWorkbook

# --- Table S2 - PERMANOVA ---
$ws2 = $wb.Worksheets.Item("Table S2 - PERMANOVA")
$ws2.Range("F3").Value = 0.09127
$ws2.Range("F4").Value = 0.00133
$ws2.Range("F9").Value = 0.23651
$ws2.Range("F12").Value = 0.74151
$ws2.Range("F14").Value = 0.00067

# --- Table S4 - Species PERMANOVA ---
$ws4 = $wb.Worksheets.Item("Table S4 - Species PERMANOVA")
$ws4.Range("F3").Value = 0.09327
$ws4.Range("F4").Value = 0.00466
$ws4.Range("F7").Value = 0.02398
$ws4.Range("F8").Value = 0.004

# --- Table S5 - HostVsymb PERMANOVA ---
$ws5 = $wb.Worksheets.Item("Table S5 - HostVsymb PERMANOVA")
$ws5.Range("F2").Value = 0.74284
$ws5.Range("F3").Value = 0.00333
$ws5.Range("K3").Value = 0.09927
$ws5.Range("F4").Value = 0.56829
$ws5.Range("F7").Value = 0.30513
$ws5.Range("K7").Value = 0.27981
$ws5.Range("F9").Value = 0.15656
$ws5.Range("F12").Value = 0.01532
$ws5.Range("K12").Value = 0.00133
$ws5.Range("F13").Value = 0.08195
$ws5.Range("F14").Value = 0.17988
$ws5.Range("K14").Value = 0.47901
